$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1446.4286
$ws.Range("I29").Value = 20.833334
$ws.Range("K29").Value = 62.500002
$ws.Range("M29").Value = 218.499998
$ws.Range("H40").Value = 1661.5385
$ws.Range("I40").Value = 1511.1111
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1511.1111
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1336.1111
$ws.Range("N40").Value = -2350
$ws.Range("H43").Value = 14870.143
$ws.Range("J43").Value = 718
$ws.Range("L43").Value = 718
$ws.Range("N43").Value = -856
$ws.Range("H62").Value = 2300
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2300
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3548
$ws.Range("H65").Value = 2300
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 11500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -17740
$ws.Range("H80").Value = 5051.087
$ws.Range("I80").Value = 178.15384
$ws.Range("K80").Value = 534.4615200000001
$ws.Range("M80").Value = 463.5384799999999
$ws.Range("H83").Value = 5051.087
$ws.Range("I83").Value = 178.15384
$ws.Range("K83").Value = 1603.38456
$ws.Range("M83").Value = 3388.61544
$ws.Range("H100").Value = 2855.4546
$ws.Range("I100").Value = 2601.111
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 2601.111
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2060.111
$ws.Range("N100").Value = -5082
$ws.Range("H132").Value = 1554.4681
$ws.Range("I132").Value = 1358.8422
$ws.Range("J132").Value = 2380.4443
$ws.Range("K132").Value = 4076.5266
$ws.Range("L132").Value = 7141.3329
$ws.Range("M132").Value = -1546.5266
$ws.Range("N132").Value = -12201.3329
$ws.Range("H141").Value = 4605.212
$ws.Range("I141").Value = 1685.5862
$ws.Range("J141").Value = 25772.5
$ws.Range("K141").Value = 5056.7586
$ws.Range("L141").Value = 77317.5
$ws.Range("M141").Value = 123.2413999999999
$ws.Range("N141").Value = -87677.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -84
$ws.Range("N4").Value = -532
$ws.Range("H74").Value = 1192.6608
$ws.Range("I74").Value = 1206.2142
$ws.Range("J74").Value = 1152
$ws.Range("K74").Value = 1206.2142
$ws.Range("L74").Value = 1152
$ws.Range("M74").Value = -332.2141999999999
$ws.Range("N74").Value = -2900
$ws.Range("H77").Value = 1192.6608
$ws.Range("I77").Value = 1206.2142
$ws.Range("J77").Value = 1152
$ws.Range("K77").Value = 6031.071
$ws.Range("L77").Value = 5760
$ws.Range("M77").Value = -1663.071
$ws.Range("N77").Value = -14496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 87210.836
$ws.Range("I86").Value = 4085.3333
$ws.Range("J86").Value = 170336.33
$ws.Range("K86").Value = 4085.3333
$ws.Range("L86").Value = 170336.33
$ws.Range("M86").Value = -2962.3333
$ws.Range("N86").Value = -172582.33
$ws.Range("H89").Value = 87210.836
$ws.Range("I89").Value = 4085.3333
$ws.Range("J89").Value = 170336.33
$ws.Range("K89").Value = 20426.6665
$ws.Range("L89").Value = 851681.6499999999
$ws.Range("M89").Value = -14810.6665
$ws.Range("N89").Value = -862913.6499999999
$ws.Range("H134").Value = 2451.8125
$ws.Range("I134").Value = 2118.32
$ws.Range("J134").Value = 3642.8572
$ws.Range("K134").Value = 6354.960000000001
$ws.Range("L134").Value = 10928.5716
$ws.Range("M134").Value = -3819.960000000001
$ws.Range("N134").Value = -15998.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1968.5441
$ws.Range("I31").Value = 2543.0334
$ws.Range("J31").Value = 1515
$ws.Range("K31").Value = 2543.0334
$ws.Range("L31").Value = 1515
$ws.Range("M31").Value = -2248.0334
$ws.Range("N31").Value = -2105
$ws.Range("H34").Value = 1968.5441
$ws.Range("I34").Value = 2543.0334
$ws.Range("J34").Value = 1515
$ws.Range("K34").Value = 2543.0334
$ws.Range("L34").Value = 1515
$ws.Range("M34").Value = -2341.0334
$ws.Range("N34").Value = -1919
$ws.Range("H58").Value = 1483271.2
$ws.Range("I58").Value = 2180452
$ws.Range("J58").Value = 1762.375
$ws.Range("K58").Value = 2180452
$ws.Range("L58").Value = 1762.375
$ws.Range("M58").Value = -2180249
$ws.Range("N58").Value = -2168.375
$ws.Range("H136").Value = 1483271.2
$ws.Range("I136").Value = 2180452
$ws.Range("J136").Value = 1762.375
$ws.Range("K136").Value = 6541356
$ws.Range("L136").Value = 5287.125
$ws.Range("M136").Value = -6538806
$ws.Range("N136").Value = -10387.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 187.25
$ws.Range("I14").Value = 187.25
$ws.Range("K14").Value = 561.75
$ws.Range("M14").Value = -388.75
$ws.Range("H17").Value = 8999.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 8999.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 26998.5
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -27336.5
$ws.Range("H68").Value = 102038.79
$ws.Range("I68").Value = 167440.75
$ws.Range("J68").Value = 1420.3846
$ws.Range("K68").Value = 502322.25
$ws.Range("L68").Value = 4261.1538
$ws.Range("M68").Value = -501511.25
$ws.Range("N68").Value = -5883.1538
$ws.Range("H71").Value = 102038.79
$ws.Range("I71").Value = 167440.75
$ws.Range("J71").Value = 1420.3846
$ws.Range("K71").Value = 1506966.75
$ws.Range("L71").Value = 12783.4614
$ws.Range("M71").Value = -1502910.75
$ws.Range("N71").Value = -20895.4614
$ws.Range("H119").Value = 7363.5454
$ws.Range("J119").Value = 10333.333
$ws.Range("L119").Value = 30999.999
$ws.Range("N119").Value = -40675.999
$ws.Range("H122").Value = 720.8461
$ws.Range("I122").Value = 614.3333
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 5528.9997
$ws.Range("L122").Value = 17991
$ws.Range("M122").Value = -3078.9997
$ws.Range("N122").Value = -22891
$ws.Range("H123").Value = 7550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1122493.9
$ws.Range("I2").Value = 2020442
$ws.Range("J2").Value = 58.75
$ws.Range("K2").Value = 2020442
$ws.Range("L2").Value = 58.75
$ws.Range("M2").Value = -2020329
$ws.Range("N2").Value = -284.75
$ws.Range("H31").Value = 23446.285
$ws.Range("I31").Value = 23446.285
$ws.Range("K31").Value = 23446.285
$ws.Range("M31").Value = -23154.285
$ws.Range("H33").Value = 16004.25
$ws.Range("I33").Value = 26008.5
$ws.Range("K33").Value = 26008.5
$ws.Range("M33").Value = -25756.5
$ws.Range("H37").Value = 23446.285
$ws.Range("I37").Value = 23446.285
$ws.Range("K37").Value = 23446.285
$ws.Range("M37").Value = -23169.285
$ws.Range("H93").Value = 31333.334
$ws.Range("J93").Value = 31333.334
$ws.Range("L93").Value = 31333.334
$ws.Range("N93").Value = -35077.334
$ws.Range("H97").Value = 34657.8
$ws.Range("I97").Value = 63891.25
$ws.Range("J97").Value = 1248.1428
$ws.Range("K97").Value = 63891.25
$ws.Range("L97").Value = 1248.1428
$ws.Range("M97").Value = -63395.25
$ws.Range("N97").Value = -2240.1428
$ws.Range("H102").Value = 3324.923
$ws.Range("I102").Value = 3323.8262
$ws.Range("J102").Value = 3333.3333
$ws.Range("K102").Value = 3323.8262
$ws.Range("L102").Value = 3333.3333
$ws.Range("M102").Value = -1701.8262
$ws.Range("N102").Value = -6577.3333
$ws.Range("H109").Value = 11707.223
$ws.Range("J109").Value = 11707.223
$ws.Range("L109").Value = 11707.223
$ws.Range("N109").Value = -13787.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3181.7144
$ws.Range("I16").Value = 711.0833
$ws.Range("J16").Value = 6475.8887
$ws.Range("K16").Value = 711.0833
$ws.Range("L16").Value = 6475.8887
$ws.Range("M16").Value = -541.0833
$ws.Range("N16").Value = -6815.8887
$ws.Range("H127").Value = 77857.5
$ws.Range("J127").Value = 77857.5
$ws.Range("L127").Value = 77857.5
$ws.Range("N127").Value = -87777.5
$ws.Range("H135").Value = 76872.11
$ws.Range("J135").Value = 76872.11
$ws.Range("L135").Value = 76872.11
$ws.Range("N135").Value = -87012.11
$ws.Range("H136").Value = 1775.7291
$ws.Range("I136").Value = 1296.8667
$ws.Range("J136").Value = 2573.8333
$ws.Range("K136").Value = 3890.6001
$ws.Range("L136").Value = 7721.499899999999
$ws.Range("M136").Value = -1340.6001
$ws.Range("N136").Value = -12821.4999
$ws.Range("H141").Value = 56666.668
$ws.Range("J141").Value = 56666.668
$ws.Range("L141").Value = 56666.668
$ws.Range("N141").Value = -67026.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22973.893
$ws.Range("J123").Value = 22973.893
$ws.Range("L123").Value = 22973.893
$ws.Range("N123").Value = -32773.893
$ws.Range("H126").Value = 8238.200000000001
$ws.Range("I126").Value = 9816.75
$ws.Range("J126").Value = 1924
$ws.Range("K126").Value = 29450.25
$ws.Range("L126").Value = 5772
$ws.Range("M126").Value = -26980.25
$ws.Range("N126").Value = -10712
